# Apply scraped schedule update (run 16:28:21) to all 3 sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item("LP1912")
$ws.Range("A2").Value = "Última actualización: 16:28:21"
$ws.Range("A3").Value = "Total filas: 255"
$ws.Range("A20").Value = "05:57:04"
$ws.Range("B20").Value = "07:21"
$ws.Range("C20").Value = "16_SANTA ANA"
$ws.Range("D20").Value = 84
$ws.Range("E20").Value = "LP1912"
$ws.Range("A21").Value = "06:16:41"
$ws.Range("B21").Value = "07:21"
$ws.Range("C21").Value = "23_HERNANDEZ"
$ws.Range("D21").Value = 65
$ws.Range("E21").Value = "LP1912"
$ws.Range("A133").Value = "10:36:50"
$ws.Range("B133").Value = "12:32"
$ws.Range("C133").Value = "14_ABASTO"
$ws.Range("D133").Value = 116
$ws.Range("E133").Value = "LP1912"
$ws.Range("A134").Value = "11:33:52"
$ws.Range("B134").Value = "12:32"
$ws.Range("C134").Value = "23_HERNANDEZ"
$ws.Range("D134").Value = 59
$ws.Range("E134").Value = "LP1912"
$ws.Range("A222").Value = "16:28:21"
$ws.Range("B222").Value = "16:32"
$ws.Range("C222").Value = "14_ABASTO"
$ws.Range("D222").Value = 4
$ws.Range("E222").Value = "LP1912"
$ws.Range("A223").Value = "15:56:56"
$ws.Range("B223").Value = "16:33"
$ws.Range("C223").Value = "83_ALUAR"
$ws.Range("D223").Value = 37
$ws.Range("E223").Value = "LP1912"
$ws.Range("A224").Value = "14:46:12"
$ws.Range("B224").Value = "16:34"
$ws.Range("C224").Value = "83_ALUAR"
$ws.Range("D224").Value = 108
$ws.Range("E224").Value = "LP1912"
$ws.Range("A225").Value = "15:56:56"
$ws.Range("B225").Value = "16:40"
$ws.Range("C225").Value = "225_GOMEZ"
$ws.Range("D225").Value = 44
$ws.Range("E225").Value = "LP1912"
$ws.Range("A226").Value = "14:46:12"
$ws.Range("B226").Value = "16:41"
$ws.Range("C226").Value = "225_GOMEZ"
$ws.Range("D226").Value = 115
$ws.Range("E226").Value = "LP1912"
$ws.Range("A227").Value = "14:53:29"
$ws.Range("B227").Value = "16:46"
$ws.Range("C227").Value = "17_ROMERO"
$ws.Range("D227").Value = 113
$ws.Range("E227").Value = "LP1912"
$ws.Range("A228").Value = "15:16:46"
$ws.Range("B228").Value = "16:53"
$ws.Range("C228").Value = "11_ETCHEVERRY"
$ws.Range("D228").Value = 97
$ws.Range("E228").Value = "LP1912"
$ws.Range("A229").Value = "16:12:06"
$ws.Range("B229").Value = "16:54"
$ws.Range("C229").Value = "11_ETCHEVERRY"
$ws.Range("D229").Value = 42
$ws.Range("E229").Value = "LP1912"
$ws.Range("A230").Value = "16:12:06"
$ws.Range("B230").Value = "16:55"
$ws.Range("C230").Value = "16_SANTA ANA"
$ws.Range("D230").Value = 43
$ws.Range("E230").Value = "LP1912"
$ws.Range("A231").Value = "15:56:56"
$ws.Range("B231").Value = "16:57"
$ws.Range("C231").Value = "15_ABASTO"
$ws.Range("D231").Value = 61
$ws.Range("E231").Value = "LP1912"
$ws.Range("A232").Value = "15:16:46"
$ws.Range("B232").Value = "16:58"
$ws.Range("C232").Value = "15_ABASTO"
$ws.Range("D232").Value = 102
$ws.Range("E232").Value = "LP1912"
$ws.Range("A233").Value = "15:56:56"
$ws.Range("B233").Value = "17:01"
$ws.Range("C233").Value = "23_HERNANDEZ"
$ws.Range("D233").Value = 65
$ws.Range("E233").Value = "LP1912"
$ws.Range("A234").Value = "15:44:42"
$ws.Range("B234").Value = "17:02"
$ws.Range("C234").Value = "23_HERNANDEZ"
$ws.Range("D234").Value = 78
$ws.Range("E234").Value = "LP1912"
$ws.Range("A235").Value = "16:28:21"
$ws.Range("B235").Value = "17:06"
$ws.Range("C235").Value = "23_HERNANDEZ"
$ws.Range("D235").Value = 38
$ws.Range("E235").Value = "LP1912"
$ws.Range("A236").Value = "15:16:46"
$ws.Range("B236").Value = "17:07"
$ws.Range("C236").Value = "16_P MOR-SANTA ANA"
$ws.Range("D236").Value = 111
$ws.Range("E236").Value = "LP1912"
$ws.Range("A237").Value = "16:28:21"
$ws.Range("B237").Value = "17:08"
$ws.Range("C237").Value = "10_OLMOS"
$ws.Range("D237").Value = 40
$ws.Range("E237").Value = "LP1912"
$ws.Range("A238").Value = "15:16:46"
$ws.Range("B238").Value = "17:09"
$ws.Range("C238").Value = "215C_EL PATO"
$ws.Range("D238").Value = 113
$ws.Range("E238").Value = "LP1912"
$ws.Range("A239").Value = "16:12:06"
$ws.Range("B239").Value = "17:10"
$ws.Range("C239").Value = "215C_EL PATO"
$ws.Range("D239").Value = 58
$ws.Range("E239").Value = "LP1912"
$ws.Range("A240").Value = "15:44:42"
$ws.Range("B240").Value = "17:21"
$ws.Range("C240").Value = "15X38_ABASTO"
$ws.Range("D240").Value = 97
$ws.Range("E240").Value = "LP1912"
$ws.Range("A241").Value = "16:28:21"
$ws.Range("B241").Value = "17:32"
$ws.Range("C241").Value = "27_EL RETIRO"
$ws.Range("D241").Value = 64
$ws.Range("E241").Value = "LP1912"
$ws.Range("A242").Value = "15:56:56"
$ws.Range("B242").Value = "17:33"
$ws.Range("C242").Value = "17_ROMERO"
$ws.Range("D242").Value = 97
$ws.Range("E242").Value = "LP1912"
$ws.Range("A243").Value = "15:44:42"
$ws.Range("B243").Value = "17:34"
$ws.Range("C243").Value = "17_ROMERO"
$ws.Range("D243").Value = 110
$ws.Range("E243").Value = "LP1912"
$ws.Range("A244").Value = "15:44:42"
$ws.Range("B244").Value = "17:36"
$ws.Range("C244").Value = "27_EL RETIRO"
$ws.Range("D244").Value = 112
$ws.Range("E244").Value = "LP1912"
$ws.Range("A245").Value = "16:12:06"
$ws.Range("B245").Value = "17:37"
$ws.Range("C245").Value = "27_EL RETIRO"
$ws.Range("D245").Value = 85
$ws.Range("E245").Value = "LP1912"
$ws.Range("A246").Value = "15:44:42"
$ws.Range("B246").Value = "17:38"
$ws.Range("C246").Value = "215B_EL PATO"
$ws.Range("D246").Value = 114
$ws.Range("E246").Value = "LP1912"
$ws.Range("A247").Value = "16:12:06"
$ws.Range("B247").Value = "17:39"
$ws.Range("C247").Value = "215B_EL PATO"
$ws.Range("D247").Value = 87
$ws.Range("E247").Value = "LP1912"
$ws.Range("A248").Value = "15:56:56"
$ws.Range("B248").Value = "17:45"
$ws.Range("C248").Value = "215_EL PELIGRO"
$ws.Range("D248").Value = 109
$ws.Range("E248").Value = "LP1912"
$ws.Range("A249").Value = "16:12:06"
$ws.Range("B249").Value = "17:46"
$ws.Range("C249").Value = "215_EL PELIGRO"
$ws.Range("D249").Value = 94
$ws.Range("E249").Value = "LP1912"
$ws.Range("A250").Value = "16:12:06"
$ws.Range("B250").Value = "17:49"
$ws.Range("C250").Value = "10_OLMOS"
$ws.Range("D250").Value = 97
$ws.Range("E250").Value = "LP1912"
$ws.Range("A251").Value = "15:56:56"
$ws.Range("B251").Value = "17:51"
$ws.Range("C251").Value = "10_OLMOS"
$ws.Range("D251").Value = 115
$ws.Range("E251").Value = "LP1912"
$ws.Range("A252").Value = "16:28:21"
$ws.Range("B252").Value = "17:52"
$ws.Range("C252").Value = "23_HERNANDEZ"
$ws.Range("D252").Value = 84
$ws.Range("E252").Value = "LP1912"
$ws.Range("A253").Value = "16:12:06"
$ws.Range("B253").Value = "17:58"
$ws.Range("C253").Value = "17_ROMERO"
$ws.Range("D253").Value = 106
$ws.Range("E253").Value = "LP1912"
$ws.Range("A254").Value = "16:28:21"
$ws.Range("B254").Value = "18:05"
$ws.Range("C254").Value = "11_ETCHEVERRY"
$ws.Range("D254").Value = 97
$ws.Range("E254").Value = "LP1912"
$ws.Range("A255").Value = "16:12:06"
$ws.Range("B255").Value = "18:06"
$ws.Range("C255").Value = "11_ETCHEVERRY"
$ws.Range("D255").Value = 114
$ws.Range("E255").Value = "LP1912"
$ws.Range("A256").Value = "16:12:06"
$ws.Range("B256").Value = "18:10"
$ws.Range("C256").Value = "16_P MOR-SANTA ANA"
$ws.Range("D256").Value = 118
$ws.Range("E256").Value = "LP1912"
$ws.Range("A257").Value = "16:12:06"
$ws.Range("B257").Value = "18:10"
$ws.Range("C257").Value = "15_ABASTO"
$ws.Range("D257").Value = 118
$ws.Range("E257").Value = "LP1912"
$ws.Range("A258").Value = "16:28:21"
$ws.Range("B258").Value = "18:17"
$ws.Range("C258").Value = "10_OLMOS"
$ws.Range("D258").Value = 109
$ws.Range("E258").Value = "LP1912"
$ws.Range("A259").Value = "16:28:21"
$ws.Range("B259").Value = "18:22"
$ws.Range("C259").Value = "215C_EL PATO"
$ws.Range("D259").Value = 114
$ws.Range("E259").Value = "LP1912"
$ws.Range("A260").Value = "16:28:21"
$ws.Range("B260").Value = "18:25"
$ws.Range("C260").Value = "16_SANTA ANA"
$ws.Range("D260").Value = 117
$ws.Range("E260").Value = "LP1912"

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Range("A2").Value = "Última actualización: 16:28:21"
$ws.Range("A3").Value = "Total filas: 43"
$ws.Range("A48").Value = "16:28:21"
$ws.Range("B48").Value = "18:22"
$ws.Range("C48").Value = "215C_EL PATO"
$ws.Range("D48").Value = 114
$ws.Range("E48").Value = "LP1912"

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Range("A2").Value = "Última actualización: 16:28:21"
$ws.Range("A3").Value = "Total filas: 39"
$ws.Range("A42").Value = "16:28:21"
$ws.Range("B42").Value = "16:31"
$ws.Range("C42").Value = "215B_LP-P MOR-40 Y 115"
$ws.Range("D42").Value = 3
$ws.Range("E42").Value = "L6173"
$ws.Range("A43").Value = "15:16:46"
$ws.Range("B43").Value = "17:05"
$ws.Range("C43").Value = "215C_LA PLATA"
$ws.Range("D43").Value = 109
$ws.Range("E43").Value = "L6203"
$ws.Range("A44").Value = "16:12:06"
$ws.Range("B44").Value = "17:06"
$ws.Range("C44").Value = "215C_LA PLATA"
$ws.Range("D44").Value = 54
$ws.Range("E44").Value = "L6203"

